# Insert two new weekly price rows for "Alcachofa" at the top of the data
# block (just above the former row 497), pushing all existing rows down by
# two. This mirrors the source diff: dimension grows from A1:R562 to
# A1:R564, and the two brand-new rows land at 497/498 with the remaining
# historical rows shifted to 499/500/.../564.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 497:562 (and everything's formatting) down by two rows, opening
# up two fresh rows at 497:498.
$ws.Rows("497:498").Insert()

# New row 497
$ws.Range("A497").Value = 9
$ws.Range("B497").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C497").Value = "Metropolitana"
$ws.Range("D497").Value = 45077
$ws.Range("E497").Value = 13
$ws.Range("F497").Value = 100112013
$ws.Range("G497").Value = "Alcachofa"
$ws.Range("H497").Value = "Española"
$ws.Range("I497").Value = "Extra"
$ws.Range("J497").Value = 43
$ws.Range("K497").Value = 19000
$ws.Range("L497").Value = 21000
$ws.Range("M497").Value = 19977
$ws.Range("N497").Value = "$/caja 25 unidades"
$ws.Range("O497").Value = "Provincia de Limarí"
$ws.Range("P497").Value = 19977
$ws.Range("Q497").Value = 1
$ws.Range("R497").Value = "Hortaliza"

# New row 498
$ws.Range("A498").Value = 9
$ws.Range("B498").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C498").Value = "Metropolitana"
$ws.Range("D498").Value = 45077
$ws.Range("E498").Value = 13
$ws.Range("F498").Value = 100112013
$ws.Range("G498").Value = "Alcachofa"
$ws.Range("H498").Value = "Española"
$ws.Range("I498").Value = "Primera"
$ws.Range("J498").Value = 70
$ws.Range("K498").Value = 18000
$ws.Range("L498").Value = 19000
$ws.Range("M498").Value = 18500
$ws.Range("N498").Value = "$/caja 30 unidades"
$ws.Range("O498").Value = "Provincia de Limarí"
$ws.Range("P498").Value = 617
$ws.Range("Q498").Value = 30
$ws.Range("R498").Value = "Hortaliza"
